$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 176, pushing the existing rows 176-180 down to 177-181.
$ws.Rows.Item(176).Insert()

$newRow = 176

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 45239
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112026
$ws.Cells.Item($newRow, 7).Value = "Haba"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 50
$ws.Cells.Item($newRow, 11).Value = 16000
$ws.Cells.Item($newRow, 12).Value = 16000
$ws.Cells.Item($newRow, 13).Value = 16000
$ws.Cells.Item($newRow, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 640
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
